$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "TestCase9" sheet by duplicating "TestCase8" (this keeps
#    column widths / styles / validations / page setup identical to its
#    sibling sheets) and then edit its contents for the new TREEVIEW test
#    case.
# ---------------------------------------------------------------------------
$src = $wb.Worksheets.Item("TestCase8")
$src.Copy($null, $src)
$newSheet = $wb.Worksheets.Item("TestCase8 (2)")
$newSheet.Name = "TestCase9"

# Relabel the TestCaseID column (was "TestCase8", now "TestCase9")
$newSheet.Range("A2:A8").Value = "TestCase9"

# Row 3 hyperlink: point at the new tree-view control documentation page
foreach ($hl in $newSheet.Hyperlinks) {
    $hl.Delete()
}
$newSheet.Range("F3").Value = "http://www.javascripttoolbox.com/lib/mktree/"
$newSheet.Hyperlinks.Add($newSheet.Range("F3"), "http://www.javascripttoolbox.com/lib/mktree/") | Out-Null

# Row 5: click a tree view item
$newSheet.Range("E5").Value = "clickTreeViewItem"
$newSheet.Range("D5").Value = "TREE_DEFAULT"
$newSheet.Range("F5").Value = "config.txt"

# Row 6: verify a tree view item exists
$newSheet.Range("F6").Value = "schedule.id"
$newSheet.Range("E6").Value = "verifyTreeViewItemExist"
$newSheet.Range("D6").Value = "TREE_DEFAULT"

# Row 7: verify a tree view item does not exist
$newSheet.Range("E7").Value = "verifyTreeViewItemNotExist"
$newSheet.Range("D7").Value = "TREE_DEFAULT"
$newSheet.Range("F7").Value = "AAAAAAAAA"

# Row 8: close the browser instead of clicking an element, and clear the old
# page-object reference that used to live in column D
$newSheet.Range("D8").ClearContents()
$newSheet.Range("E8").Value = "closeBrowser"

# The old sheet had two extra rows (radio-button checks) that don't apply to
# the tree-view test case, so drop them.
$newSheet.Rows.Item(10).Delete()
$newSheet.Rows.Item(9).Delete()

$newSheet.Range("C13").Select()

# ---------------------------------------------------------------------------
# 2. "TestCase8" sheet gains one more step: close the browser at the end of
#    the test (TS_010).
# ---------------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item("TestCase8")
$ws8.Range("A11").Value = "TestCase8"
$ws8.Range("B11").Value = "TS_010"
$ws8.Range("E11").Value = "closeBrowser"
$ws8.Range("D10").Copy() | Out-Null
$ws8.Range("F2").Copy()
$ws8.Range("D10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws8.Range("E11").Select()

# ---------------------------------------------------------------------------
# 3. "TestCase7" sheet: just a cursor/selection move, no data change.
# ---------------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("TestCase7")
$ws7.Range("E7").Select()

# ---------------------------------------------------------------------------
# 4. "TestSuite" sheet: TestCase9 is now flagged to run ("Yes" instead of
#    "No" in the Run Mode column).
# ---------------------------------------------------------------------------
$wsSuite = $wb.Worksheets.Item("TestSuite")
$wsSuite.Range("C10").Value = "Yes"
$wsSuite.Range("C8").Select()

$wsSuite.Activate()
